$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 tweaks ---
$ws.Range("AD2").Value = 8
$ws.Range("AO2").Value = 8

# --- Row 4 becomes the match previously on row 6 (with a few odds re-priced) ---
$ws.Range("A4").Value = "IPYYCLH5"
$ws.Range("B4").Value = "28/11/2024"
$ws.Range("C4").Value = "20:30"
$ws.Range("D4").Value = "PARAGUAY - PRIMERA DIVISION"
$ws.Range("E4").Value = "Cerro Porteno"
$ws.Range("F4").Value = "Libertad Asuncion"
$ws.Range("G4").Value = 1.83
$ws.Range("H4").Value = 3.3
$ws.Range("I4").Value = 4.1
$ws.Range("J4").Value = 2.6
$ws.Range("K4").Value = 2.05
$ws.Range("L4").Value = 4.75
$ws.Range("M4").Value = 1.07
$ws.Range("N4").Value = 9
$ws.Range("O4").Value = 1.36
$ws.Range("P4").Value = 3
$ws.Range("Q4").Value = 2.2
$ws.Range("R4").Value = 1.65
$ws.Range("S4").Value = 1.44
$ws.Range("T4").Value = 2.63
$ws.Range("U4").Value = 2
$ws.Range("V4").Value = 1.73
$ws.Range("W4").Value = 6
$ws.Range("X4").Value = 8
$ws.Range("Y4").Value = 9
$ws.Range("Z4").Value = 15
$ws.Range("AA4").Value = 17
$ws.Range("AB4").Value = 34
$ws.Range("AC4").Value = 8.5
$ws.Range("AD4").Value = 6.5
$ws.Range("AE4").Value = 17
$ws.Range("AF4").Value = 67
$ws.Range("AG4").Value = 1000
$ws.Range("AH4").Value = 10
$ws.Range("AI4").Value = 21
$ws.Range("AJ4").Value = 15
$ws.Range("AK4").Value = 41
$ws.Range("AL4").Value = 41
$ws.Range("AM4").Value = 41
$ws.Range("AN4").Value = 3.75
$ws.Range("AO4").Value = 10
$ws.Range("AP4").Value = 23
$ws.Range("AQ4").Value = 34
$ws.Range("AR4").Value = 51
$ws.Range("AS4").Value = 201
$ws.Range("AT4").Value = 2.63
$ws.Range("AU4").Value = 9
$ws.Range("AV4").Value = 67
$ws.Range("AW4").Value = 6
$ws.Range("AX4").Value = 23
$ws.Range("AY4").Value = 34
$ws.Range("AZ4").Value = 81
$ws.Range("BA4").Value = 126
$ws.Range("BB4").Value = 301

# --- Rows 5 and 6 (the two fixtures now folded away) are removed entirely ---
$ws.Rows("5:6").Delete()
